$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) mirroring the structure of row 3.

$ws.Range("A4").Value = 42633.676747685182
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9948
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.12
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -1.04
$ws.Range("I4").Value = $false

# Match the date-formatted style already used for A3/G3 (style index 1)
# instead of letting a new NumberFormat assignment create a duplicate style.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
